$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 1798
$ws.Range("I29").Value = 90
$ws.Range("J29").Value = 2225
$ws.Range("K29").Value = 270
$ws.Range("L29").Value = 6675
$ws.Range("M29").Value = 11
$ws.Range("N29").Value = -7237
# Row 32
$ws.Range("H32").Value = 588.1177
$ws.Range("I32").Value = 400
$ws.Range("J32").Value = 666.5
$ws.Range("K32").Value = 400
$ws.Range("L32").Value = 666.5
$ws.Range("M32").Value = -74
$ws.Range("N32").Value = -1318.5
# Row 40
$ws.Range("H40").Value = 3038170.2
$ws.Range("I40").Value = 5084559.5
$ws.Range("J40").Value = 1430292.9
$ws.Range("K40").Value = 5084559.5
$ws.Range("L40").Value = 1430292.9
$ws.Range("M40").Value = -5084384.5
$ws.Range("N40").Value = -1430642.9
# Row 62
$ws.Range("H62").Value = 2182.5
$ws.Range("I62").Value = 2089
$ws.Range("J62").Value = 2650
$ws.Range("K62").Value = 2089
$ws.Range("L62").Value = 2650
$ws.Range("M62").Value = -1465
$ws.Range("N62").Value = -3898
# Row 65
$ws.Range("H65").Value = 2182.5
$ws.Range("I65").Value = 2089
$ws.Range("J65").Value = 2650
$ws.Range("K65").Value = 10445
$ws.Range("L65").Value = 13250
$ws.Range("M65").Value = -7325
$ws.Range("N65").Value = -19490
# Row 137
$ws.Range("H137").Value = 2542.8572
$ws.Range("I137").Value = 2475
$ws.Range("K137").Value = 7425
$ws.Range("M137").Value = -4875
# Row 138
$ws.Range("H138").Value = 1633.3549
$ws.Range("I138").Value = 783.36365
$ws.Range("J138").Value = 2100.85
$ws.Range("K138").Value = 2350.09095
$ws.Range("L138").Value = 6302.549999999999
$ws.Range("M138").Value = 2789.90905
$ws.Range("N138").Value = -16582.55

$ws = $wb.Worksheets.Item("ARM")
# Row 132
$ws.Range("H132").Value = 7637.9443
$ws.Range("I132").Value = 9122.416999999999
$ws.Range("K132").Value = 27367.251
$ws.Range("M132").Value = -24837.251
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").ClearContents()
$ws.Range("N133").Value = 0

$ws = $wb.Worksheets.Item("BSM")
# Row 63
$ws.Range("H63").Value = 44700
$ws.Range("I63").Value = 49400
$ws.Range("J63").Value = 40000
$ws.Range("K63").Value = 49400
$ws.Range("L63").Value = 40000
$ws.Range("M63").Value = -48714
$ws.Range("N63").Value = -41372
# Row 66
$ws.Range("H66").Value = 44700
$ws.Range("I66").Value = 49400
$ws.Range("J66").Value = 40000
$ws.Range("K66").Value = 148200
$ws.Range("L66").Value = 120000
$ws.Range("M66").Value = -144768
$ws.Range("N66").Value = -126864
# Row 107
$ws.Range("H107").Value = 850.44446
$ws.Range("I107").Value = 629.3570999999999
$ws.Range("J107").Value = 1624.25
$ws.Range("K107").Value = 629.3570999999999
$ws.Range("L107").Value = 1624.25
$ws.Range("M107").Value = 1290.6429
$ws.Range("N107").Value = -5464.25
# Row 126
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("N126").Value = 0

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 9527628
$ws.Range("I31").Value = 3167.5789
$ws.Range("J31").Value = 100010000
$ws.Range("K31").Value = 3167.5789
$ws.Range("L31").Value = 100010000
$ws.Range("M31").Value = -2872.5789
$ws.Range("N31").Value = -100010590
# Row 34
$ws.Range("H34").Value = 9527628
$ws.Range("I34").Value = 3167.5789
$ws.Range("J34").Value = 100010000
$ws.Range("K34").Value = 3167.5789
$ws.Range("L34").Value = 100010000
$ws.Range("M34").Value = -2965.5789
$ws.Range("N34").Value = -100010404
# Row 44
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").ClearContents()
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = 0

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 1595.5
$ws.Range("J34").Value = 1724.9286
$ws.Range("L34").Value = 5174.7858
$ws.Range("N34").Value = -5342.7858
# Row 87
$ws.Range("H87").Value = 1000
$ws.Range("I87").Value = 1000
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 3000
$ws.Range("L87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -1752
# Row 90
$ws.Range("H90").Value = 1000
$ws.Range("I90").Value = 1000
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 9000
$ws.Range("L90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -2760

$ws = $wb.Worksheets.Item("GSM")
# Row 53
$ws.Range("H53").Value = 9633.333000000001
$ws.Range("J53").Value = 9633.333000000001
$ws.Range("L53").Value = 9633.333000000001
$ws.Range("N53").Value = -10895.333
# Row 55
$ws.Range("H55").Value = 6666.6665
$ws.Range("I55").Value = 10000
$ws.Range("J55").Value = 5000
$ws.Range("K55").Value = 10000
$ws.Range("L55").Value = 5000
$ws.Range("M55").Value = -9673
$ws.Range("N55").Value = -5654
# Row 62
$ws.Range("H62").Value = 50000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 50000
$ws.Range("K62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("M62").Value = 50000
$ws.Range("N62").Value = -51372
# Row 65
$ws.Range("H65").Value = 50000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 50000
$ws.Range("K65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("M65").Value = 150000
$ws.Range("N65").Value = -156864
# Row 80
$ws.Range("H80").Value = 6121.8887
$ws.Range("I80").Value = 4360.4
$ws.Range("K80").Value = 4360.4
$ws.Range("M80").Value = -3362.4
# Row 83
$ws.Range("H83").Value = 6121.8887
$ws.Range("I83").Value = 4360.4
$ws.Range("K83").Value = 21802
$ws.Range("M83").Value = -16810
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").ClearContents()
$ws.Range("N123").Value = 0
# Row 125
$ws.Range("H125").Value = 36666.668
$ws.Range("J125").Value = 36666.668
$ws.Range("L125").Value = 36666.668
$ws.Range("N125").Value = -41586.668
# Row 126
$ws.Range("H126").Value = 7874
$ws.Range("I126").Value = 6998.4
$ws.Range("J126").Value = 9333.333000000001
$ws.Range("K126").Value = 20995.2
$ws.Range("L126").Value = 27999.999
$ws.Range("M126").Value = -18525.2
$ws.Range("N126").Value = -32939.999
# Row 127
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").ClearContents()
$ws.Range("N127").Value = 0
# Row 128
$ws.Range("H128").Value = 150000
$ws.Range("J128").Value = 50000
$ws.Range("L128").Value = 50000
$ws.Range("N128").Value = -59960
# Row 129
$ws.Range("H129").Value = 49399.4
$ws.Range("J129").Value = 49399.4
$ws.Range("L129").Value = 49399.4
$ws.Range("N129").Value = -59399.4

$ws = $wb.Worksheets.Item("LTW")
# Row 70
$ws.Range("H70").Value = 24500
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 24500
$ws.Range("K70").Value = 0
$ws.Range("L70").ClearContents()
$ws.Range("M70").Value = 24500
$ws.Range("N70").Value = -25040
# Row 73
$ws.Range("H73").Value = 24500
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 24500
$ws.Range("K73").Value = 0
$ws.Range("L73").ClearContents()
$ws.Range("M73").Value = 24500
$ws.Range("N73").Value = -26372
# Row 87
$ws.Range("H87").Value = 39666.668
$ws.Range("J87").Value = 39666.668
$ws.Range("L87").Value = 39666.668
$ws.Range("N87").Value = -41912.668
# Row 90
$ws.Range("H90").Value = 39666.668
$ws.Range("J90").Value = 39666.668
$ws.Range("L90").Value = 119000.004
$ws.Range("N90").Value = -130232.004
# Row 93
$ws.Range("H93").Value = 1690494
$ws.Range("I93").Value = 2080447
$ws.Range("J93").Value = 698
$ws.Range("K93").Value = 2080447
$ws.Range("L93").Value = 698
$ws.Range("M93").Value = -2079199
$ws.Range("N93").Value = -3194
# Row 122
$ws.Range("H122").Value = 2367.75
$ws.Range("I122").Value = 2333.6667
$ws.Range("J122").Value = 2470
$ws.Range("K122").Value = 7001.000100000001
$ws.Range("L122").Value = 7410
$ws.Range("M122").Value = -4551.000100000001
$ws.Range("N122").Value = -12310

$ws = $wb.Worksheets.Item("WVR")
# Row 111
$ws.Range("H111").Value = 58583.715
$ws.Range("J111").Value = 58583.715
$ws.Range("L111").Value = 58583.715
$ws.Range("N111").Value = -66763.715
# Row 126
$ws.Range("H126").Value = 1368.421
$ws.Range("I126").Value = 1390
$ws.Range("J126").Value = 1344.4445
$ws.Range("K126").Value = 4170
$ws.Range("L126").Value = 4033.3335
$ws.Range("M126").Value = -1700
$ws.Range("N126").Value = -8973.333500000001
# Row 138
$ws.Range("H138").Value = 81742.664
$ws.Range("J138").Value = 81742.664
$ws.Range("L138").Value = 81742.664
$ws.Range("N138").Value = -92022.664
